$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows 367-374 (update through 9/09/2021 inclusive)
$data = @(
    @{ Row = 367; A = 44441; B = 2; C = 5; D = 32.29348317509527 },
    @{ Row = 368; A = 44442; B = 0; C = 4; D = 25.83478654007622 },
    @{ Row = 369; A = 44443; B = 4; C = 6; D = 38.75217981011431 },
    @{ Row = 370; A = 44444; B = 1; C = 7; D = 45.21087644513337 },
    @{ Row = 371; A = 44445; B = 1; C = 8; D = 51.66957308015243 },
    @{ Row = 372; A = 44446; B = 1; C = 9; D = 58.12826971517148 },
    @{ Row = 373; A = 44447; B = 0; C = 9; D = 58.12826971517148 },
    @{ Row = 374; A = 44448; B = 0; C = 7; D = 45.21087644513337 }
)

foreach ($item in $data) {
    $r = $item.Row

    # Copy the date-column formatting (style s="2") from the row above
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}
